# "Generate Report for Handoff"
#
# The 530290e6-5491-4005-a69b-99d51fd2293c.md file moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with a fresh
# handoff timestamp recorded, and the localization report grew an
# "Error Detail" note (plus a wider column to show it) on the per-language
# sheets explaining that the handback file seen isn't the latest.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66600597a83b1112b9500398b9017861a869914f/e2e/530290e6-5491-4005-a69b-99d51fd2293c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fa81ab0dec71bd6709a6023685a61080103bfa21/e2e/530290e6-5491-4005-a69b-99d51fd2293c.md."

# --- Overview sheet: row 3 is the 530290e6-...-md file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusNew
$overview.Range("F3").Value = $statusNew
$overview.Range("G3").Value = "2016-09-06 06:56:44"

# Excel's ColumnWidth (characters) round-trips to the OOXML <col width="..">
# (points-ish units) with a fixed +5/6 padding offset for this workbook's
# font, so back it out to land on an even 40 in the saved file.
$col16Width = 40 - (5 / 6)

# --- zh-cn sheet: row 3 is the 530290e6-...-md file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusNew
$zhcn.Range("H3").Value = "2016-09-06 06:56:39"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = $col16Width

# --- de-de sheet: row 3 is the 530290e6-...-md file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusNew
$dede.Range("H3").Value = "2016-09-06 06:56:44"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = $col16Width
